# Append 9 new master-data rows (regcntr_id 10002-10010 / machine_id
# 10021-10029) below the existing table, mirroring the existing row layout
# (lang_code = "eng", is_active = TRUE, cr_by = "superadmin",
# cr_dtimes = eff_dtimes = "now()"), then leave the sheet selection on the
# newly-entered machine_id column, matching how the data was pasted in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$regcntrStart = 10002
$machineStart = 10021
$firstRow = 22
$lastRow = 30

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $offset = $row - $firstRow

    $ws.Cells.Item($row, 1).Value = $regcntrStart + $offset
    $ws.Cells.Item($row, 2).Value = $machineStart + $offset
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

$ws.Range("B" + $firstRow + ":B" + $lastRow).Select() | Out-Null

# Page setup was touched to portrait orientation (as printed in the source
# workbook's page-setup metadata).
$ws.PageSetup.Orientation = 1
